$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.089.76"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.655.18"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D5").Value = "'217.63"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'0.5267"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.2610"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.06347"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "'20.42"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "'0.07798"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "'4.506"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "1.674.86"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "'0.5495"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "0.0₅8219"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "'65.41"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "26.119.73"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'190.70"
$ws.Range("D21").Value = "'10.06"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'141.56"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").Value = "'0.1235"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'7.238"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "'16.06"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'0.05883"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").Value = "'1.273"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "'3.522"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "'0.9504"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D36").Value = "'2.411"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").Value = "'0.5707"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "'0.01615"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("D39").Value = "'5.796"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").Value = "'0.8490"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "1.027.01"
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("D44").Value = "1.799.45"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'57.08"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "'1.0000"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").Value = "'1.475"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").Value = "'7.847"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("D50").Value = "'0.05152"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").Value = "'0.09693"
$ws.Range("E51").Value = "  -0.30%  "
